$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Find-RowByLabel($sheet, $label) {
    $lastRow = $sheet.Cells.Item($sheet.Rows.Count, 1).End(-4162).Row
    for ($r = 1; $r -le $lastRow; $r++) {
        if ($sheet.Cells.Item($r, 1).Value2 -eq $label) {
            return $r
        }
    }
    return -1
}

# Delete entire rows for "SC 92" then "RM 232" (delete higher row first so
# the other row's index doesn't shift before we look it up / delete it).
$rowSC92 = Find-RowByLabel $ws "SC 92"
if ($rowSC92 -gt 0) {
    $ws.Rows.Item($rowSC92).Delete() | Out-Null
}

$rowRM232 = Find-RowByLabel $ws "RM 232"
if ($rowRM232 -gt 0) {
    $ws.Rows.Item($rowRM232).Delete() | Out-Null
}

# Apply individual cell edits (values now at their final row numbers since
# both deleted rows were below row 25, so rows 1-25 are unaffected).
$ws.Range("F5").Value = $null
$ws.Range("F11").Value = 17.65
$ws.Range("D19").Value = -15.5
$ws.Range("F19").Value = $null
$ws.Range("D21").Value = $null
$ws.Range("D23").Value = -13.9
$ws.Range("F25").Value = 16.6
$ws.Range("D27").Value = $null
$ws.Range("F29").Value = $null
$ws.Range("D33").Value = -14.1
